# Updates "3-) Turkiye'de tekstil ithalat ve ihracat hacmi.xlsx"
# - Splits the old "Türkiyenin Toplam Ticaret Hacmi" column into two columns:
#     "Türkiyenin Toplam İhracat Hacmi" (reuses column B, values replaced)
#     "Türkiyenin Toplam İthalat Hacmi" (new column inserted at D)
# - Shifts old D/E/F (Tekstil İthalat Hacmi / oranlari) one column right to E/F/G
# - Recalculates the ratio columns (F, G) against the new export/import totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; it inherits the formatting of the old D.
$ws.Columns("D").Insert()

# ---- Header row ----
$ws.Cells.Item(1, 2).Value2 = "Türkiyenin Toplam İhracat Hacmi"
$ws.Cells.Item(1, 4).Value2 = "Türkiyenin Toplam İthalat Hacmi"

# ---- Data rows ----
# Columns: A Yıl | B Toplam İhracat Hacmi | C Tekstil İhracat Hacmi |
#          D Toplam İthalat Hacmi | E Tekstil İthalat Hacmi |
#          F Tekstil İhracat Oranı (%) | G Tekstil İthalat Oranı (%)
$rows = @(
    @{ Row = 2;  B = 160179713.31900001;  D = 210252947.71900001;  F = 18.78870221790449;  G = 6.6406173000057684 }
    @{ Row = 3;  B = 169897259.215;       D = 234868399.14500001;  F = 19.736116670703531; G = 6.7308362148967884 }
    @{ Row = 4;  B = 176114435.222;       D = 228745844.17899999;  F = 20.09000099531886;  G = 7.2672155958346876 }
    @{ Row = 5;  B = 161303390.877;       D = 200121709.28400001;  F = 19.968089050626809; G = 7.0291230178517461 }
    @{ Row = 6;  B = 162033589.39500001;  D = 196478323.14899999;  F = 19.927004715848351; G = 6.9098854318418983 }
    @{ Row = 7;  B = 170238045.59599999;  D = 222559771.331;       F = 18.85367236837812;  G = 6.3563664948956244 }
    @{ Row = 8;  B = 188343441.39899999;  D = 210232429.727;       F = 18.24143356721229;  G = 6.2571493936886986 }
    @{ Row = 9;  B = 190669538.46900001;  D = 194382889.35100001;  F = 18.442103422155739; G = 6.4801034695264947 }
    @{ Row = 10; B = 179376774.61500001;  D = 209080024.24200001;  F = 18.57366408137765;  G = 5.1425387561439422 }
    @{ Row = 11; B = 233216566.11500001;  D = 248305934.92399999;  F = 17.896233588920659; G = 5.4912996643327876 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value2 = $r.B   # Türkiyenin Toplam İhracat Hacmi
    # Column C (Tekstil İhracat Hacmi) is unchanged.
    $ws.Cells.Item($row, 4).Value2 = $r.D   # Türkiyenin Toplam İthalat Hacmi (new)
    # Column E (old column D, Tekstil İthalat Hacmi) is unchanged; shifted by the insert.
    $ws.Cells.Item($row, 6).Value2 = $r.F   # Tekstil İhracat Oranı (%), recomputed
    $ws.Cells.Item($row, 7).Value2 = $r.G   # Tekstil İthalat Oranı (%), recomputed
}

# ---- Column widths (bestFit widths recorded after the edit) ----
# Target stored widths (OOXML "width" attr): D=30.140625 E=19.140625 F=22.42578125 G=22.140625
# The interop column-width setter quantizes to whole pixels, so we pick the nearest
# achievable input that rounds to those stored widths.
$ws.Columns("D").ColumnWidth = 29.3
$ws.Columns("E").ColumnWidth = 18.3
$ws.Columns("F").ColumnWidth = 21.65
$ws.Columns("G").ColumnWidth = 21.3

# ---- Selection state matches the file as last saved ----
$ws.Range("F12").Select()
